# Update the worker rank matrix (female) sheet: refresh the computed
# "matrices" score (column E) for every worker, which in turn re-orders the
# ranking. Rows keep their sheet position (A/G are just positional
# index/rank counters) but the person (name/prolific id/race) occupying
# each row changes to reflect the new sort order by column E (desc).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, prolificId(B), name(C), score(E), race(F)
$rows = @(
  @{ Row = 2;  B = 2;  C = "Colleen";   E = 13.41179440177212;  F = "White" },
  @{ Row = 3;  B = 3;  C = "Annes";     E = 13.28345529631224;  F = "Asian" },
  @{ Row = 4;  B = 19; C = "Jewel";     E = 8.23763560135623;   F = "Black or African American" },
  @{ Row = 5;  B = 22; C = "Khushi";    E = 8.166612813012238;  F = "Asian" },
  @{ Row = 6;  B = 21; C = "Bri";       E = 8.005597717419404;  F = "Black or African American" },
  @{ Row = 7;  B = 30; C = "Shadaisia"; E = 5.361927025870195;  F = "Black or African American" },
  @{ Row = 8;  B = 32; C = "Kellie";    E = 5.237683582040133;  F = "White" },
  @{ Row = 9;  B = 33; C = "Shaniek";   E = 5.010962683506764;  F = "Black or African American" },
  @{ Row = 10; B = 35; C = "Lori";      E = 4.477479880056773;  F = "White" },
  @{ Row = 11; B = 34; C = "Tina";      E = 4.389849861394186;  F = "White" },
  @{ Row = 12; B = 41; C = "Giana";     E = 2.382527936458554;  F = "White" },
  @{ Row = 13; B = 44; C = "Nansi";     E = 1.029173221199296;  F = "Asian" }
)

foreach ($r in $rows) {
  $ws.Cells.Item($r.Row, 2).Value = $r.B
  $ws.Cells.Item($r.Row, 3).Value = $r.C
  $ws.Cells.Item($r.Row, 5).Value = $r.E
  $ws.Cells.Item($r.Row, 6).Value = $r.F
}
